$wb = $excel.ActiveWorkbook

# --- "Btts" sheet: append a new match prediction as row 19 ---
$btts = $wb.Worksheets.Item("Btts")
$btts.Cells.Item(19, 1).Value = "16-12-2024 20:45"
$btts.Cells.Item(19, 2).Value = "PORTUGAL"
$btts.Cells.Item(19, 3).Value = "PRIMEIRA LIGA"
$btts.Cells.Item(19, 4).Value = "Rio Ave - Guimaraes"
$btts.Cells.Item(19, 5).Value = 76.7
$btts.Cells.Item(19, 6).Value = 1.95

# --- "Over_Under" sheet: remove the stale Honduras match (row 5), ---
# --- the rows below it shift up automatically.                    ---
$overUnder = $wb.Worksheets.Item("Over_Under")
$overUnder.Rows("5:5").Delete()
